$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.538.04'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -4.47%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.808.12'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -3.20%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '274.73'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -8.54%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5003'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -5.93%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3416'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -8.35%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '44.00'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.99%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06623'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -7.56%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.44'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -9.27%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7961'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -10.17%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07833'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.808.65'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.16%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.009'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.25%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '86.31'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -6.40%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.001'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.09%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.93'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -5.90%  '
$ws.Range('E19').Value = '  +0.11%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007929'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -6.46%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '25.595.64'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.39%  '
$ws.Range('E22').Value = '  -5.43%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.845'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -7.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.074'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.50%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.240'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.16%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '142.34'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.32%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.658'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.17%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.98'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.71%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '108.14'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.72%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.243'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -9.39%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.185'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -9.34%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08692'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.63%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04765'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.122'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.88%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.851'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.12%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7091'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -11.38%  '
$ws.Range('B37').Value = 'MXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.110'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.03%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.304'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -13.42%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01822'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -6.02%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.4999'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -17.79%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9302'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -12.42%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '115.41'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.64%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.133'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -5.32%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.001'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.736'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -11.18%  '
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1341'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -9.99%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4343'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -16.34%  '
$ws.Range('B48').Value = 'Elrond'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '35.98'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.66%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.129'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -8.04%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05816'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.99%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.463'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -10.41%  '
